$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4525806
$ws.Range("C3").Value = 1955395
$ws.Range("C4").Value = 274396
$ws.Range("C10").Value = 1872350
$ws.Range("C11").Value = 732510
$ws.Range("C12").Value = 197264
